# Update crypto price/volume data per commit "Updated cryptos list on Thu Aug 24 08:45:36 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.487.50"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").Value = "1.674.56"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'220.29"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").Value = "'0.5312"
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +3.06%  "
$ws.Range("D9").Value = "'0.06389"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").Value = "'21.73"
$ws.Range("E10").Value = "  +4.22%  "
$ws.Range("D11").Value = "'0.07797"
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").Value = "1.680.49"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "0.0₅8338"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "'65.63"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "26.507.45"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "'4.760"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").Value = "'192.72"
$ws.Range("E20").Value = "  +2.34%  "
$ws.Range("D21").Value = "'10.33"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").Value = "'6.314"
$ws.Range("E22").Value = "  +1.36%  "
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'0.1274"
$ws.Range("E24").Value = "  +4.46%  "
$ws.Range("D25").Value = "'138.05"
$ws.Range("E25").Value = "  -5.54%  "
$ws.Range("D26").Value = "'7.409"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  +2.85%  "
$ws.Range("D28").Value = "'1.428"
$ws.Range("E28").Value = "  +3.48%  "
$ws.Range("D29").Value = "'0.06260"
$ws.Range("E29").Value = "  +4.49%  "
$ws.Range("D30").Value = "'1.291"
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("E31").Value = "  +5.80%  "
$ws.Range("D32").Value = "'3.425"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").Value = "'1.693"
$ws.Range("E33").Value = "  +1.81%  "
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("D35").Value = "'0.6156"
$ws.Range("E35").Value = "  +8.93%  "
$ws.Range("D36").Value = "'2.419"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("D37").Value = "'2.781"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("D38").Value = "'6.133"
$ws.Range("E38").Value = "  +4.86%  "
$ws.Range("D39").Value = "'0.01618"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "1.092.62"
$ws.Range("E40").Value = "  +6.06%  "
$ws.Range("D41").Value = "'0.8613"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "'100.70"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").Value = "1.820.28"
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'58.76"
$ws.Range("E45").Value = "  +5.20%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'8.195"
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.004"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.516"
$ws.Range("E48").Value = "  +9.31%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05195"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'6.017"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.4233"
$ws.Range("E51").Value = "  +0.26%  "
